# Update the "Classroom" column (Q) on the "Results" sheet with new
# classroom assignments (0-3) for a number of students, reflecting a
# rebalanced/more dynamic assignment that now also accounts for
# additional characteristics beyond gender/special needs/grades.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Results")

# Map of row number -> new Classroom (column Q) value
$changes = @{
    2 = 1
    3 = 2
    4 = 0
    5 = 3
    7 = 2
    8 = 3
    9 = 1
    10 = 2
    12 = 3
    13 = 1
    14 = 2
    15 = 3
    16 = 2
    19 = 3
    21 = 3
    23 = 1
    24 = 2
    25 = 0
    26 = 3
    28 = 3
    29 = 2
    30 = 0
    31 = 1
    33 = 2
    34 = 2
    35 = 3
    36 = 0
    37 = 1
    38 = 3
    39 = 3
    40 = 0
    41 = 2
    42 = 0
    43 = 1
    44 = 3
    45 = 1
    46 = 1
    47 = 1
    51 = 2
    52 = 2
    53 = 2
    54 = 3
    55 = 1
    56 = 0
    57 = 0
    58 = 3
    59 = 1
    61 = 0
    62 = 2
    63 = 1
    64 = 3
    66 = 3
    69 = 0
    70 = 1
    71 = 3
    73 = 3
    74 = 1
    75 = 2
    77 = 3
    79 = 2
    80 = 3
    81 = 1
    82 = 0
    84 = 1
    87 = 3
    88 = 2
    89 = 1
    90 = 2
    91 = 3
    92 = 2
    93 = 0
    94 = 1
    95 = 1
    97 = 0
    99 = 1
    100 = 2
}

foreach ($row in $changes.Keys) {
    $ws.Cells.Item([int]$row, 17).Value = $changes[$row]
}

Write-Output "Updated $($changes.Count) Classroom (Q) values on Results sheet"
